# Manage Interviewers - Bug Fixes
# Updates the "AMS" sheet (Graph Sprint Data / Quick Interview History Data):
#  - corrects the B12 run-time timestamp
#  - applies explicit "Normal" styling to A12/C12/E12/F12/G12 (previously unstyled)
#  - fills in the previously-blank rows 13 and 14 with new interview-history
#    data rows (2021-06-17 runs of "145_data_hstry" / "145_hstry_data")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "AMS" sheet

# Helper-less approach: write date-like / text values through a literal-text
# formula and then collapse it to a plain value via copy/paste-special so that
# Excel does not auto-convert the "2021-06-17" style text into a date serial.

# ---------------------------------------------------------------------------
# Row 12 - existing row: correct the B12 timestamp and make the formatting
# explicit (matching the "Normal" style used elsewhere on the sheet).
# ---------------------------------------------------------------------------
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = 44363.80313471065
$ws.Range("C12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Style = "Normal"

# ---------------------------------------------------------------------------
# Row 13 - new data row (styled like row 11)
# ---------------------------------------------------------------------------
$ws.Range("A13").Formula = "=""2021-06-17"""
$ws.Range("A13").Copy()
$ws.Range("A13").PasteSpecial(-4163)
$ws.Range("A13").Style = "Normal"

$ws.Range("B13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B13").Value = 44364.57115336806

$ws.Range("C13").Formula = "=""145_data_hstry"""
$ws.Range("C13").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$ws.Range("C13").Style = "Normal"

$ws.Range("D13").Value = 96

$ws.Range("E13").Style = "Normal"
$ws.Range("E13").Value = 95

$ws.Range("F13").Style = "Normal"
$ws.Range("F13").Value = 1

$ws.Range("G13").Style = "Normal"
$ws.Range("G13").Value = 2.48

# ---------------------------------------------------------------------------
# Row 14 - new data row (left unstyled, like row 12 originally was)
# ---------------------------------------------------------------------------
$ws.Range("A14").Formula = "=""2021-06-17"""
$ws.Range("A14").Copy()
$ws.Range("A14").PasteSpecial(-4163)

$ws.Range("B14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B14").Value = 44364.61015453223

$ws.Range("C14").Formula = "=""145_hstry_data"""
$ws.Range("C14").Copy()
$ws.Range("C14").PasteSpecial(-4163)

$ws.Range("D14").Value = 96
$ws.Range("E14").Value = 95
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.51

$excel.CutCopyMode = 0
